$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"4.269024333333333"
$ws.Range("H2").Value = [double]"12.807073"
$ws.Range("I2").Value = [double]"0.239419143716247"
$ws.Range("J2").Value = [double]"0.2624375843396372"
$ws.Range("M2").Value = [double]"1.009860666666667"
$ws.Range("N2").Value = [double]"3.029582"
$ws.Range("O2").Value = [double]"0.01353413605720072"
$ws.Range("P2").Value = [double]"0.01542521070970148"
$ws.Range("Q2").Value = [double]"4.311119759276223"
$ws.Range("R2").Value = [double]"38.800077833486"
$ws.Range("S2").Value = [double]"0.00324033126575418"
$ws.Range("T2").Value = [double]"0.004048155036583958"
$ws.Range("G3").Value = [double]"4.269024333333333"
$ws.Range("H3").Value = [double]"12.807073"
$ws.Range("I3").Value = [double]"0.239419143716247"
$ws.Range("J3").Value = [double]"0.2624375843396372"
$ws.Range("O3").Value = [double]"0.6185519418990597"
$ws.Range("P3").Value = [double]"0.704979911415303"
$ws.Range("Q3").Value = [double]"197.0315273608429"
$ws.Range("R3").Value = [double]"1773.283746247586"
$ws.Range("S3").Value = [double]"0.1480931762734946"
$ws.Range("T3").Value = [double]"0.1850132249598035"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = [double]"4.269024333333333"
$ws.Range("H4").Value = [double]"12.807073"
$ws.Range("I4").Value = [double]"0.239419143716247"
$ws.Range("J4").Value = [double]"0.2624375843396372"
$ws.Range("K4").Value = [double]"1"
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.009315666666666667"
$ws.Range("N4").Value = [double]"0.027947"
$ws.Range("O4").Value = [double]"0.0001248484115599408"
$ws.Range("P4").Value = [double]"0.000142293017222847"
$ws.Range("Q4").Value = [double]"0.03976880768122222"
$ws.Range("R4").Value = [double]"0.357919269131"
$ws.Range("S4").Value = [double]"2.989109979001462E-05"
$ws.Range("T4").Value = [double]"3.734303570836236E-05"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = [double]"4.269024333333333"
$ws.Range("H5").Value = [double]"12.807073"
$ws.Range("I5").Value = [double]"0.239419143716247"
$ws.Range("J5").Value = [double]"0.2624375843396372"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"27.4428835"
$ws.Range("N5").Value = [double]"54.885767"
$ws.Range("O5").Value = [double]"0.3677890736321797"
$ws.Range("P5").Value = [double]"0.2794525848577725"
$ws.Range("Q5").Value = [double]"117.1543374383318"
$ws.Range("R5").Value = [double]"702.926024629991"
$ws.Range("S5").Value = [double]"0.08805574507720816"
$ws.Range("T5").Value = [double]"0.0733388613075413"
$ws.Range("I6").Value = [double]"0.1495657278917394"
$ws.Range("J6").Value = [double]"0.1639454043592587"
$ws.Range("M6").Value = [double]"1.009860666666667"
$ws.Range("N6").Value = [double]"3.029582"
$ws.Range("O6").Value = [double]"0.01353413605720072"
$ws.Range("P6").Value = [double]"0.01542521070970148"
$ws.Range("Q6").Value = [double]"2.693167116113333"
$ws.Range("R6").Value = [double]"24.23850404502"
$ws.Range("S6").Value = [double]"0.002024242910781062"
$ws.Range("T6").Value = [double]"0.002528892407128777"
$ws.Range("I7").Value = [double]"0.1495657278917394"
$ws.Range("J7").Value = [double]"0.1639454043592587"
$ws.Range("O7").Value = [double]"0.6185519418990597"
$ws.Range("P7").Value = [double]"0.704979911415303"
$ws.Range("S7").Value = [double]"0.09251417142898177"
$ws.Range("T7").Value = [double]"0.1155782166421362"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = [double]"0.1495657278917394"
$ws.Range("J8").Value = [double]"0.1639454043592587"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.009315666666666667"
$ws.Range("N8").Value = [double]"0.027947"
$ws.Range("O8").Value = [double]"0.0001248484115599408"
$ws.Range("P8").Value = [double]"0.000142293017222847"
$ws.Range("Q8").Value = [double]"0.02484367196333333"
$ws.Range("R8").Value = [double]"0.22359304767"
$ws.Range("S8").Value = [double]"1.867304355109E-05"
$ws.Range("T8").Value = [double]"2.332828624609862E-05"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = [double]"0.1495657278917394"
$ws.Range("J9").Value = [double]"0.1639454043592587"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"27.4428835"
$ws.Range("N9").Value = [double]"54.885767"
$ws.Range("O9").Value = [double]"0.3677890736321797"
$ws.Range("P9").Value = [double]"0.2794525848577725"
$ws.Range("Q9").Value = [double]"73.18660271964499"
$ws.Range("R9").Value = [double]"439.11961631787"
$ws.Range("S9").Value = [double]"0.05500864050842549"
$ws.Range("T9").Value = [double]"0.04581496702374758"
$ws.Range("G10").Value = [double]"3.607962333333333"
$ws.Range("H10").Value = [double]"10.823887"
$ws.Range("I10").Value = [double]"0.2023448884238746"
$ws.Range("J10").Value = [double]"0.221798904202795"
$ws.Range("M10").Value = [double]"1.009860666666667"
$ws.Range("N10").Value = [double]"3.029582"
$ws.Range("O10").Value = [double]"0.01353413605720072"
$ws.Range("P10").Value = [double]"0.01542521070970148"
$ws.Range("Q10").Value = [double]"3.643539247248222"
$ws.Range("R10").Value = [double]"32.791853225234"
$ws.Range("S10").Value = [double]"0.002738563250407818"
$ws.Range("T10").Value = [double]"0.003421294832509006"
$ws.Range("G11").Value = [double]"3.607962333333333"
$ws.Range("H11").Value = [double]"10.823887"
$ws.Range("I11").Value = [double]"0.2023448884238746"
$ws.Range("J11").Value = [double]"0.221798904202795"
$ws.Range("O11").Value = [double]"0.6185519418990597"
$ws.Range("P11").Value = [double]"0.704979911415303"
$ws.Range("Q11").Value = [double]"166.5210300270149"
$ws.Range("R11").Value = [double]"1498.689270243134"
$ws.Range("S11").Value = [double]"0.1251608236679362"
$ws.Range("T11").Value = [double]"0.1563637718368977"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = [double]"3.607962333333333"
$ws.Range("H12").Value = [double]"10.823887"
$ws.Range("I12").Value = [double]"0.2023448884238746"
$ws.Range("J12").Value = [double]"0.221798904202795"
$ws.Range("K12").Value = [double]"1"
$ws.Range("L12").Value = [double]"0.3333333333333333"
$ws.Range("M12").Value = [double]"0.009315666666666667"
$ws.Range("N12").Value = [double]"0.027947"
$ws.Range("O12").Value = [double]"0.0001248484115599408"
$ws.Range("P12").Value = [double]"0.000142293017222847"
$ws.Range("Q12").Value = [double]"0.03361057444322222"
$ws.Range("R12").Value = [double]"0.3024951699889999"
$ws.Range("S12").Value = [double]"2.52624379069942E-05"
$ws.Range("T12").Value = [double]"3.15604352957369E-05"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = [double]"3.607962333333333"
$ws.Range("H13").Value = [double]"10.823887"
$ws.Range("I13").Value = [double]"0.2023448884238746"
$ws.Range("J13").Value = [double]"0.221798904202795"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"27.4428835"
$ws.Range("N13").Value = [double]"54.885767"
$ws.Range("O13").Value = [double]"0.3677890736321797"
$ws.Range("P13").Value = [double]"0.2794525848577725"
$ws.Range("Q13").Value = [double]"99.01288998605483"
$ws.Range("R13").Value = [double]"594.0773399163289"
$ws.Range("S13").Value = [double]"0.0744202390676236"
$ws.Range("T13").Value = [double]"0.06198227709809253"
$ws.Range("G14").Value = [double]"4.691815"
$ws.Range("H14").Value = [double]"9.38363"
$ws.Range("I14").Value = [double]"0.2631304584056895"
$ws.Range("J14").Value = [double]"0.1922857150526861"
$ws.Range("M14").Value = [double]"1.009860666666667"
$ws.Range("N14").Value = [double]"3.029582"
$ws.Range("O14").Value = [double]"0.01353413605720072"
$ws.Range("P14").Value = [double]"0.01542521070970148"
$ws.Range("Q14").Value = [double]"4.738079423776667"
$ws.Range("R14").Value = [double]"28.42847654266"
$ws.Range("S14").Value = [double]"0.003561243424856198"
$ws.Range("T14").Value = [double]"0.002966047671153301"
$ws.Range("G15").Value = [double]"4.691815"
$ws.Range("H15").Value = [double]"9.38363"
$ws.Range("I15").Value = [double]"0.2631304584056895"
$ws.Range("J15").Value = [double]"0.1922857150526861"
$ws.Range("O15").Value = [double]"0.6185519418990597"
$ws.Range("P15").Value = [double]"0.704979911415303"
$ws.Range("Q15").Value = [double]"216.5449066022767"
$ws.Range("R15").Value = [double]"1299.26943961366"
$ws.Range("S15").Value = [double]"0.162759856019629"
$ws.Range("T15").Value = [double]"0.1355575663642709"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = [double]"4.691815"
$ws.Range("H16").Value = [double]"9.38363"
$ws.Range("I16").Value = [double]"0.2631304584056895"
$ws.Range("J16").Value = [double]"0.1922857150526861"
$ws.Range("K16").Value = [double]"1"
$ws.Range("L16").Value = [double]"0.3333333333333333"
$ws.Range("M16").Value = [double]"0.009315666666666667"
$ws.Range("N16").Value = [double]"0.027947"
$ws.Range("O16").Value = [double]"0.0001248484115599408"
$ws.Range("P16").Value = [double]"0.000142293017222847"
$ws.Range("Q16").Value = [double]"0.04370738460166667"
$ws.Range("R16").Value = [double]"0.26224430761"
$ws.Range("S16").Value = [double]"3.285141976498941E-05"
$ws.Range("T16").Value = [double]"2.736091456369931E-05"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = [double]"4.691815"
$ws.Range("H17").Value = [double]"9.38363"
$ws.Range("I17").Value = [double]"0.2631304584056895"
$ws.Range("J17").Value = [double]"0.1922857150526861"
$ws.Range("K17").Value = [double]"2"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"27.4428835"
$ws.Range("N17").Value = [double]"54.885767"
$ws.Range("O17").Value = [double]"0.3677890736321797"
$ws.Range("P17").Value = [double]"0.2794525848577725"
$ws.Range("Q17").Value = [double]"128.7569324485525"
$ws.Range("R17").Value = [double]"515.02772979421"
$ws.Range("S17").Value = [double]"0.09677650754143934"
$ws.Range("T17").Value = [double]"0.05373474010269823"
$ws.Range("G18").Value = [double]"2.595084333333333"
$ws.Range("H18").Value = [double]"7.785253"
$ws.Range("I18").Value = [double]"0.1455397815624493"
$ws.Range("J18").Value = [double]"0.159532392045623"
$ws.Range("M18").Value = [double]"1.009860666666667"
$ws.Range("N18").Value = [double]"3.029582"
$ws.Range("O18").Value = [double]"0.01353413605720072"
$ws.Range("P18").Value = [double]"0.01542521070970148"
$ws.Range("Q18").Value = [double]"2.620673594916223"
$ws.Range("R18").Value = [double]"23.586062354246"
$ws.Range("S18").Value = [double]"0.001969755205401462"
$ws.Range("T18").Value = [double]"0.00246082076232644"
$ws.Range("G19").Value = [double]"2.595084333333333"
$ws.Range("H19").Value = [double]"7.785253"
$ws.Range("I19").Value = [double]"0.1455397815624493"
$ws.Range("J19").Value = [double]"0.159532392045623"
$ws.Range("O19").Value = [double]"0.6185519418990597"
$ws.Range("P19").Value = [double]"0.704979911415303"
$ws.Range("Q19").Value = [double]"119.7729012304829"
$ws.Range("R19").Value = [double]"1077.956111074346"
$ws.Range("S19").Value = [double]"0.090023914509018"
$ws.Range("T19").Value = [double]"0.1124671316121947"
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = [double]"2.595084333333333"
$ws.Range("H20").Value = [double]"7.785253"
$ws.Range("I20").Value = [double]"0.1455397815624493"
$ws.Range("J20").Value = [double]"0.159532392045623"
$ws.Range("K20").Value = [double]"1"
$ws.Range("L20").Value = [double]"0.3333333333333333"
$ws.Range("M20").Value = [double]"0.009315666666666667"
$ws.Range("N20").Value = [double]"0.027947"
$ws.Range("O20").Value = [double]"0.0001248484115599408"
$ws.Range("P20").Value = [double]"0.000142293017222847"
$ws.Range("Q20").Value = [double]"0.02417494062122222"
$ws.Range("R20").Value = [double]"0.217574465591"
$ws.Range("S20").Value = [double]"1.817041054685256E-05"
$ws.Range("T20").Value = [double]"2.270034540894982E-05"
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = [double]"2.595084333333333"
$ws.Range("H21").Value = [double]"7.785253"
$ws.Range("I21").Value = [double]"0.1455397815624493"
$ws.Range("J21").Value = [double]"0.159532392045623"
$ws.Range("K21").Value = [double]"2"
$ws.Range("L21").Value = [double]"1"
$ws.Range("M21").Value = [double]"27.4428835"
$ws.Range("N21").Value = [double]"54.885767"
$ws.Range("O21").Value = [double]"0.3677890736321797"
$ws.Range("P21").Value = [double]"0.2794525848577725"
$ws.Range("Q21").Value = [double]"71.21659703234184"
$ws.Range("R21").Value = [double]"427.299582194051"
$ws.Range("S21").Value = [double]"0.05352794143748302"
$ws.Range("T21").Value = [double]"0.04458173932569291"
